$d = $word.ActiveDocument

$d.Content.Find.Execute("Answers: Logarithms", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Answers: Logarithms", 2)

$d.Content.Find.Execute("Zoë Gemmell", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Zoë Gemmell", 2)

$d.Content.Find.Execute("Answers to questions relating to the study guide on logarithms.", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Answers to questions relating to the study guide on logarithms.", 2)
